$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab09")

$ws.Range("G62").Value = 48.279545454545399

$ws.Range("C63").Value = 2.3285714285714301
$ws.Range("D63").Value = 8.5481481481481492
$ws.Range("E63").Value = 22.8603174603175
$ws.Range("F63").Value = 34.976190476190503
$ws.Range("G63").Value = 42.574603174603098

$ws.Range("C65").Value = 5.0128205128205101
$ws.Range("D65").Value = 21.815384615384598
$ws.Range("E65").Value = 51.815384615384602
$ws.Range("F65").Value = 35.105128205128203
$ws.Range("G65").Value = 43.210256410256399

$ws.Range("C66").Value = 10.0677419354839
$ws.Range("D66").Value = 21.931899641577001
$ws.Range("E66").Value = 39.736917562724102
$ws.Range("F66").Value = 37.116487455197102
$ws.Range("G66").Value = 44.387003610108302

$ws.Range("G68").Value = 45.038095238095302

$ws.Range("C76").Value = 2.1615384615384601
$ws.Range("D76").Value = 13.015384615384599
$ws.Range("E76").Value = 39.769230769230802
$ws.Range("F76").Value = 37.269230769230802
$ws.Range("G76").Value = 44.815384615384602

$ws.Range("D82").Value = 51.130769230769303

$ws.Range("C83").Value = 1.97167630057803
$ws.Range("D83").Value = 7.6254335260115704
$ws.Range("E83").Value = 21.5549132947977
$ws.Range("F83").Value = 35.036416184971102
$ws.Range("G83").Value = 42.587283236994203

$ws.Range("E86").Value = 72.535000000000096

$ws.Range("C87").Value = 5.4872340425531902
$ws.Range("D87").Value = 22.463829787234101
$ws.Range("E87").Value = 54.731914893617002
$ws.Range("F87").Value = 35.865957446808501
$ws.Range("G87").Value = 43.638297872340502

$ws.Range("F89").Value = 38.133333333333397

$ws.Range("E91").Value = 85.757142857142796

$ws.Range("D95").Value = 68.1933333333333
